$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EJ45FF")

# Making EJ45FF the active sheet/tab (mirrors activeTab 1 -> 4, tabSelected
# moving off EJ44 and onto EJ45FF).
$ws.Activate()

# Header row: rename the generic "String N" headers to the EJ45FF-specific
# part numbers (new shared strings J4501FF..J4506FF).
$ws.Range("B1").Value = "J4501FF"
$ws.Range("C1").Value = "J4502FF"
$ws.Range("D1").Value = "J4503FF"
$ws.Range("E1").Value = "J4504FF"
$ws.Range("F1").Value = "J4505FF"
$ws.Range("G1").Value = "J4506FF"

# Fill in the previously-empty measurement data.
$ws.Range("B2").Value = 291.7
$ws.Range("C2").Value = 215
$ws.Range("D2").Value = 174
$ws.Range("E2").Value = 131.9
$ws.Range("F2").Value = 98.2
$ws.Range("G2").Value = 73.8

$ws.Range("B3").Value = 297.5
$ws.Range("C3").Value = 220.4
$ws.Range("D3").Value = 179.1
$ws.Range("E3").Value = 134.8
$ws.Range("F3").Value = 100.5
$ws.Range("G3").Value = 76

$ws.Range("B4").Value = 304.4
$ws.Range("C4").Value = 226.7
$ws.Range("D4").Value = 183.5
$ws.Range("E4").Value = 138.1
$ws.Range("F4").Value = 102.7
$ws.Range("G4").Value = 77.6

$ws.Range("B5").Value = 310.5
$ws.Range("C5").Value = 232.3
$ws.Range("D5").Value = 188.5
$ws.Range("E5").Value = 141.2
$ws.Range("F5").Value = 105
$ws.Range("G5").Value = 79.2

$ws.Range("B6").Value = 315.9
$ws.Range("C6").Value = 236.1
$ws.Range("D6").Value = 192.5
$ws.Range("E6").Value = 143.9
$ws.Range("F6").Value = 107.8
$ws.Range("G6").Value = 80.8

$ws.Range("B7").Value = 320.7
$ws.Range("C7").Value = 241.8
$ws.Range("D7").Value = 195.8
$ws.Range("E7").Value = 146.5
$ws.Range("F7").Value = 110.2
$ws.Range("G7").Value = 81.9

$ws.Range("B8").Value = 327.1
$ws.Range("C8").Value = 246.5
$ws.Range("D8").Value = 199.9
$ws.Range("E8").Value = 149.3
$ws.Range("F8").Value = 112.5
$ws.Range("G8").Value = 83.5

# Final selection on the now-active EJ45FF tab.
[void]$ws.Range("G9").Select()
